# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.9
$ws.Range("G2").Value = 3.35
$ws.Range("H2").Value = 2.46
$ws.Range("I2").Value = 2.82
$ws.Range("K2").Value = 3.6
$ws.Range("P2").Value = 1.72
$ws.Range("BH2").Value = "'2026-02-21 03:42:55"

# Row 3
$ws.Range("BH3").Value = "'2026-02-21 03:42:55"

# Row 4
$ws.Range("H4").Value = 1.56
$ws.Range("I4").Value = 1.64
$ws.Range("K4").Value = 4.6
$ws.Range("M4").Value = 1.07
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 1.92
$ws.Range("Q4").Value = 1.91
$ws.Range("U4").Value = 1.84
$ws.Range("V4").Value = 2.56
$ws.Range("X4").Value = 17
$ws.Range("Z4").Value = 9.2
$ws.Range("AC4").Value = 9.8
$ws.Range("AK4").Value = 120
$ws.Range("AN4").Value = 160
$ws.Range("AO4").Value = 10.5
$ws.Range("AX4").Value = 7.8
$ws.Range("AY4").Value = 19
$ws.Range("BA4").Value = 7.4
$ws.Range("BB4").Value = 8.4
$ws.Range("BC4").Value = 8.4
$ws.Range("BD4").Value = 8.4
$ws.Range("BE4").Value = 8.4
$ws.Range("BF4").Value = 8.4
$ws.Range("BG4").Value = 8
$ws.Range("BH4").Value = "'2026-02-21 03:42:55"

# Row 5
$ws.Range("F5").Value = 3.6
$ws.Range("G5").Value = 4.8
$ws.Range("H5").Value = 1.96
$ws.Range("I5").Value = 2.18
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 4.8
$ws.Range("P5").Value = 1.95
$ws.Range("Q5").Value = 1.64
$ws.Range("BH5").Value = "'2026-02-21 03:42:55"

# Row 6
$ws.Range("BA6").Value = 14
$ws.Range("BH6").Value = "'2026-02-21 03:42:55"

# Row 7
$ws.Range("G7").Value = 2.16
$ws.Range("P7").Value = 2.08
$ws.Range("Q7").Value = 1.81
$ws.Range("BH7").Value = "'2026-02-21 03:42:55"

# Row 8
$ws.Range("P8").Value = 1.97
$ws.Range("Q8").Value = 1.79
$ws.Range("BH8").Value = "'2026-02-21 03:42:55"

# Row 9
$ws.Range("P9").Value = 2.88
$ws.Range("BH9").Value = "'2026-02-21 03:42:55"

# Row 10
$ws.Range("BH10").Value = "'2026-02-21 03:42:55"

# Row 11
$ws.Range("BH11").Value = "'2026-02-21 03:42:55"

# Row 12
$ws.Range("BH12").Value = "'2026-02-21 03:42:55"

# Row 13
$ws.Range("BH13").Value = "'2026-02-21 03:42:55"

# Row 14
$ws.Range("P14").Value = 2.44
$ws.Range("Q14").Value = 1.51
$ws.Range("BH14").Value = "'2026-02-21 03:42:55"

# Row 15
$ws.Range("BH15").Value = "'2026-02-21 03:42:55"

# Row 16
$ws.Range("BH16").Value = "'2026-02-21 03:42:55"

# Row 17
$ws.Range("BH17").Value = "'2026-02-21 03:42:55"

# Row 18
$ws.Range("BH18").Value = "'2026-02-21 03:42:55"

# Row 19
$ws.Range("G19").Value = 3.05
$ws.Range("P19").Value = 1.73
$ws.Range("BH19").Value = "'2026-02-21 03:42:55"

# Row 20
$ws.Range("Q20").Value = 2.2
$ws.Range("U20").Value = 1.91
$ws.Range("Y20").Value = 14.5
$ws.Range("AS20").Value = 12.5
$ws.Range("BE20").Value = 13
$ws.Range("BG20").Value = 12.5
$ws.Range("BH20").Value = "'2026-02-21 03:42:55"

# Row 21
$ws.Range("F21").Value = 3.95
$ws.Range("H21").Value = 1.98
$ws.Range("I21").Value = 2
$ws.Range("N21").Value = 4.6
$ws.Range("Q21").Value = 1.76
$ws.Range("AI21").Value = 29
$ws.Range("AT21").Value = 16
$ws.Range("BE21").Value = 50
$ws.Range("BH21").Value = "'2026-02-21 03:42:55"

# Row 22
$ws.Range("A22").Value = "'Portuguese Primeira Liga"
$ws.Range("C22").Value = "'17:15:00"
$ws.Range("D22").Value = "'Famalicao"
$ws.Range("E22").Value = "'Casa Pia"
$ws.Range("F22").Value = 1.32
$ws.Range("G22").Value = 1.75
$ws.Range("H22").Value = 5.7
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 3.85
$ws.Range("K22").Value = 980
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 1.25
$ws.Range("Q22").Value = 1.01
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("T22").Value = 0
$ws.Range("U22").Value = 0
$ws.Range("X22").Value = 0
$ws.Range("Y22").Value = 0
$ws.Range("Z22").Value = 0
$ws.Range("AA22").Value = 0
$ws.Range("AB22").Value = 0
$ws.Range("AC22").Value = 0
$ws.Range("AD22").Value = 0
$ws.Range("AE22").Value = 0
$ws.Range("AF22").Value = 0
$ws.Range("AG22").Value = 0
$ws.Range("AH22").Value = 0
$ws.Range("AI22").Value = 0
$ws.Range("AJ22").Value = 0
$ws.Range("AK22").Value = 0
$ws.Range("AL22").Value = 0
$ws.Range("AM22").Value = 0
$ws.Range("AN22").Value = 0
$ws.Range("AO22").Value = 0
$ws.Range("AP22").Value = 0
$ws.Range("AQ22").Value = 0
$ws.Range("AR22").Value = 0
$ws.Range("AS22").Value = 0
$ws.Range("AT22").Value = 0
$ws.Range("AU22").Value = 0
$ws.Range("AV22").Value = 0
$ws.Range("AW22").Value = 0
$ws.Range("AX22").Value = 0
$ws.Range("AY22").Value = 0
$ws.Range("AZ22").Value = 0
$ws.Range("BA22").Value = 0
$ws.Range("BB22").Value = 0
$ws.Range("BC22").Value = 0
$ws.Range("BD22").Value = 0
$ws.Range("BE22").Value = 0
$ws.Range("BF22").Value = 0
$ws.Range("BG22").Value = 0
$ws.Range("BH22").Value = "'2026-02-21 03:42:55"

# Row 23
$ws.Range("A23").Value = "'Portuguese Segunda Liga"
$ws.Range("D23").Value = "'Porto B"
$ws.Range("E23").Value = "'Pacos Ferreira"
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("P23").Value = 1.24
$ws.Range("BH23").Value = "'2026-02-21 03:42:55"

# Row 24
$ws.Range("A24").Value = "'Chilean Primera B"
$ws.Range("C24").Value = "'18:00:00"
$ws.Range("D24").Value = "'Cobreloa Calama"
$ws.Range("E24").Value = "'Deportes Temuco"
$ws.Range("F24").Value = 1.58
$ws.Range("G24").Value = 2.02
$ws.Range("H24").Value = 1.98
$ws.Range("I24").Value = 1000
$ws.Range("J24").Value = 1.98
$ws.Range("K24").Value = 1000
$ws.Range("P24").Value = 1.6
$ws.Range("Q24").Value = 1.94
$ws.Range("BH24").Value = "'2026-02-21 03:42:55"

# Row 25
$ws.Range("A25").Value = "'Colombian Primera A"
$ws.Range("D25").Value = "'Cucuta Deportivo"
$ws.Range("E25").Value = "'Tolima"
$ws.Range("F25").Value = 2.42
$ws.Range("G25").Value = 3.15
$ws.Range("H25").Value = 2.62
$ws.Range("I25").Value = 3.6
$ws.Range("J25").Value = 2.96
$ws.Range("K25").Value = 5.1
$ws.Range("P25").Value = 1.48
$ws.Range("Q25").Value = 2.22
$ws.Range("BH25").Value = "'2026-02-21 03:42:55"

# Row 26
$ws.Range("A26").Value = "'Chilean Primera B"
$ws.Range("C26").Value = "'20:30:00"
$ws.Range("D26").Value = "'Curico Unido"
$ws.Range("E26").Value = "'Magallanes"
$ws.Range("F26").Value = 1.04
$ws.Range("G26").Value = 1000
$ws.Range("H26").Value = 1.04
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 1.01
$ws.Range("K26").Value = 1000
$ws.Range("P26").Value = 1.24
$ws.Range("Q26").Value = 2.22
$ws.Range("BH26").Value = "'2026-02-21 03:42:55"

# Delete trailing row 27 (event removed from source feed)
$ws.Rows(27).Delete()

Write-Host "edit complete"